# Generate Report for Handoff
# The handoff run picked up a new source GUID/commit and re-stamped the
# handoff/handback timestamps. Mirror that into the report workbook.

$wb = $excel.ActiveWorkbook

$newGuid = "67ca1668-05ae-4bbd-a9ae-6c043f570a87"

$newZhXlf = "$newGuid.bed7c293627cf97d4932911678118c2d963683d3.zh-cn.xlf"
$newDeXlf = "$newGuid.bed7c293627cf97d4932911678118c2d963683d3.de-de.xlf"

$newFileName = "$newGuid.md"
$newPathName = "e2e\$newGuid.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# A2 (File Name) is plain text.
$wsOverview.Range("A2").Value = $newFileName
# B2 (Path And Name) carries the e2e\<guid>.md hyperlink display text.
$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newPathName
# G2: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-25 22:56:09"

# --- zh-cn sheet ------------------------------------------------------
# A2 (Source File Name) carries the <guid>.md hyperlink display text.
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = $newFileName
# G2: Latest Handoff File
$wsZhCn.Range("G2").Value = $newZhXlf
# H2: Latest Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-08-25 22:55:59"

# --- de-de sheet ------------------------------------------------------
# A2 (Source File Name) carries the <guid>.md hyperlink display text.
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = $newFileName
# G2: Latest Handoff File
$wsDeDe.Range("G2").Value = $newDeXlf
# H2: Latest Handoff Datetime (shares the Overview G2 timestamp)
$wsDeDe.Range("H2").Value = "2016-08-25 22:56:09"
